# Add two new columns, I ("I0") and J ("IF"), to the results table on Sheet1.
# Header cells pick up the same formatting (bold font, thin border, centered
# alignment) already used by the other header cells (copied from H1), and the
# data rows (2-40) get the corresponding numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column data for rows 2-40 (I = "I0", J = "IF").
$iValues = @(7,4,6,7,7,5,7,5,4,5,6,7,10,3,7,8,5,6,6,7,8,9,7,5,9,9,6,9,9,7,7,7,7,6,8,8,5,8,7)
$jValues = @(7,4,7,7,7,5,7,5,5,6,6,7,10,4,7,8,5,7,7,7,8,9,7,5,9,9,6,9,9,7,7,7,7,6,8,8,5,8,7)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
